# Insert a new data row at row 161 (shifts existing rows 161-218 down to 162-219)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("161:161").Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A161").Value = 11
$ws.Range("B161").Value = "Vega Monumental Concepción"
$ws.Range("C161").Value = "Bíobío"
$ws.Range("D161").Value = 45093
$ws.Range("E161").Value = 8
$ws.Range("F161").Value = 100112043
$ws.Range("G161").Value = "Pepino ensalada"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 100
$ws.Range("K161").Value = 12000
$ws.Range("L161").Value = 13000
$ws.Range("M161").Value = 12500
$ws.Range("N161").Value = "$/caja 60 unidades"
$ws.Range("O161").Value = "Región de Arica y Parinacota"
$ws.Range("P161").Value = 208
$ws.Range("Q161").Value = 60
$ws.Range("R161").Value = "Hortaliza"
